# Apply the betexplorer "chile primera-division 2023" update:
#  - three pairs of adjacent rows (45/46, 67/68, 96/97) have their match
#    details (columns F:V) swapped, while the shared match date/time in
#    columns A:E stays put;
#  - two brand-new match rows are appended at the end of the sheet
#    (rows 189 and 190), extending the used range from A1:V188 to A1:V190.
#
# NOTE: this runtime's PowerShell-style function calls only bind
# parameters positionally (named "-param value" binding does not work),
# so every helper below is invoked with plain positional arguments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchDetails($rowA, $rowB) {
    # Columns F..V hold the match details (teams, scores, odds, timestamps,
    # url); columns A..E (index/pais/torneio/temporada/data_partida) are
    # identical for both rows in every swapped pair, so they are untouched.
    $cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
    foreach ($col in $cols) {
        $refA = "$col$rowA"
        $refB = "$col$rowB"
        $valA = $ws.Range($refA).Value2
        $valB = $ws.Range($refB).Value2
        $ws.Range($refA).Value2 = $valB
        $ws.Range($refB).Value2 = $valA
    }
}

Swap-MatchDetails 45 46
Swap-MatchDetails 67 68
Swap-MatchDetails 96 97

# Writes one full data row, matching the layout used by every other row in
# the sheet:
#   row, indice, dataPartida,
#   home, homeGols, away, awayGols,
#   homeOpenOdds, homeOpenDate, homeCloseOdds, homeCloseDate,
#   drawOpenOdds, drawOpenDate, drawCloseOdds, drawCloseDate,
#   awayOpenOdds, awayOpenDate, awayCloseOdds, awayCloseDate,
#   url
function Set-MatchRow {
    param(
        $row, $indice, $dataPartida,
        $home, $homeGols, $away, $awayGols,
        $homeOpenOdds, $homeOpenDate, $homeCloseOdds, $homeCloseDate,
        $drawOpenOdds, $drawOpenDate, $drawCloseOdds, $drawCloseDate,
        $awayOpenOdds, $awayOpenDate, $awayCloseOdds, $awayCloseDate,
        $url
    )

    # Copy formatting (border/bold on A, date format on E, ...) from the
    # previous row so the new row looks like the rest of the table.
    $prev = $row - 1
    $ws.Range("A${prev}:V${prev}").Copy()
    $ws.Range("A${row}:V${row}").PasteSpecial(-4122)

    $ws.Range("A$row").Value2 = $indice
    $ws.Range("B$row").Value2 = "chile"
    $ws.Range("C$row").Value2 = "primera-division"

    # "temporada" must stay the text value "2023", not get auto-converted
    # to a number by the object model. Copying the value straight from the
    # previous row (which already holds it as text) avoids that coercion
    # without touching number formats / styles at all.
    $ws.Range("D${prev}").Copy()
    $ws.Range("D$row").PasteSpecial(-4163)

    $ws.Range("E$row").Value2 = $dataPartida

    $ws.Range("F$row").Value2 = $home
    $ws.Range("G$row").Value2 = $homeGols
    $ws.Range("H$row").Value2 = $away
    $ws.Range("I$row").Value2 = $awayGols

    $ws.Range("J$row").Value2 = $homeOpenOdds
    $ws.Range("K$row").Value2 = $homeOpenDate
    $ws.Range("L$row").Value2 = $homeCloseOdds
    $ws.Range("M$row").Value2 = $homeCloseDate

    $ws.Range("N$row").Value2 = $drawOpenOdds
    $ws.Range("O$row").Value2 = $drawOpenDate
    $ws.Range("P$row").Value2 = $drawCloseOdds
    $ws.Range("Q$row").Value2 = $drawCloseDate

    $ws.Range("R$row").Value2 = $awayOpenOdds
    $ws.Range("S$row").Value2 = $awayOpenDate
    $ws.Range("T$row").Value2 = $awayCloseOdds
    $ws.Range("U$row").Value2 = $awayCloseDate

    $ws.Range("V$row").Value2 = $url
}

Set-MatchRow 189 188 45193.9375 `
    "U. Catolica" 1 "Magallanes" 0 `
    1.73 "16/09/2023 21:42" 2.32 "24/09/2023 22:27" `
    3.88 "16/09/2023 21:42" 3.36 "24/09/2023 22:30" `
    4.34 "16/09/2023 21:42" 3.28 "24/09/2023 22:30" `
    "https://www.betexplorer.com/football/chile/primera-division/u-catolica-magallanes/pA1X4IAr/"

Set-MatchRow 190 189 45194.04166666666 `
    "A. Italiano" 2 "Union La Calera" 1 `
    2.39 "17/09/2023 00:44" 2.22 "25/09/2023 00:53" `
    3.36 "17/09/2023 00:44" 3.41 "25/09/2023 00:53" `
    2.9 "17/09/2023 00:44" 3.46 "25/09/2023 00:53" `
    "https://www.betexplorer.com/football/chile/primera-division/a-italiano-union-la-calera/lvG1AvmL/"
